$d = $word.ActiveDocument

$find = "Datas das campanhas de Constelação de Pégaso 2022: 8 a 17 de outubro, 7 a 16 de novembro,"
$replace = "Datas das campanhas de 2022 que usam Constelação de Pégaso: 8 a 17 de outubro, 7 a 16 de novembro,"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
